$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (subject IDs) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated meanEMG values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 3.7733248942804378
$ws.Range("D2").Value = 0.15343052869178794
$ws.Range("E2").Value = 3.592424634487986

# Row 3 (STR) updated meanEMG values
$ws.Range("B3").Value = 0.51057939332039604
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 0.44673890480000134
$ws.Range("E3").Value = 6.1638797954227664

# Update selection to match the new active range
$ws.Range("B1:E3").Select()
